# Update symbol list (cryptos.xlsx) with refreshed Price (D) and
# Volume(1h) (E) figures, as pulled on Fri Feb 3 06:58:12 UTC 2023.
#
# Values are written with a leading apostrophe so Excel stores them
# as literal text (matching the existing text cells) instead of
# reinterpreting them as numbers, which would otherwise round or
# reformat values such as "8.650" or "0.00000000751". The style is
# reset to Normal afterwards so the quote-prefix formatting does not
# linger on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "323.07"
Set-TextValue $ws.Range("E2") "-2.06%"

Set-TextValue $ws.Range("D3") "39.57"
Set-TextValue $ws.Range("E3") "-1.35%"

Set-TextValue $ws.Range("D4") "5.868"
Set-TextValue $ws.Range("E4") "11.29%"

Set-TextValue $ws.Range("D5") "0.08025"
Set-TextValue $ws.Range("E5") "-0.78%"

Set-TextValue $ws.Range("D6") "8.646"

Set-TextValue $ws.Range("D7") "1.987"
Set-TextValue $ws.Range("E7") "2.49%"

Set-TextValue $ws.Range("D8") "2.942"
Set-TextValue $ws.Range("E8") "-0.52%"

Set-TextValue $ws.Range("D9") "0.9281"
Set-TextValue $ws.Range("E9") "-0.88%"

Set-TextValue $ws.Range("D10") "0.1281"
Set-TextValue $ws.Range("E10") "-6.29%"

Set-TextValue $ws.Range("D11") "0.1962"
Set-TextValue $ws.Range("E11") "-0.64%"

Set-TextValue $ws.Range("D12") "8.661"
Set-TextValue $ws.Range("E12") "33.41%"

Set-TextValue $ws.Range("D13") "0.09169"
Set-TextValue $ws.Range("E13") "0.84%"

Set-TextValue $ws.Range("D14") "0.03562"
Set-TextValue $ws.Range("E14") "1.62%"

Set-TextValue $ws.Range("E15") "9.07%"

Set-TextValue $ws.Range("D16") "0.001291"
Set-TextValue $ws.Range("E16") "-7.84%"

Set-TextValue $ws.Range("D17") "0.006270"
Set-TextValue $ws.Range("E17") "-1.33%"

Set-TextValue $ws.Range("D18") "3.350"
Set-TextValue $ws.Range("E18") "-0.47%"

Set-TextValue $ws.Range("D19") "4.578"
Set-TextValue $ws.Range("E19") "1.19%"

Set-TextValue $ws.Range("E20") "0.55%"

Set-TextValue $ws.Range("E21") "4.35%"

Set-TextValue $ws.Range("E22") "-4.77%"

Set-TextValue $ws.Range("D23") "0.04406"
Set-TextValue $ws.Range("E23") "-0.98%"

Set-TextValue $ws.Range("E24") "3.33%"

Set-TextValue $ws.Range("D25") "0.004394"
Set-TextValue $ws.Range("E25") "1.62%"

Set-TextValue $ws.Range("E26") "-11.68%"

Set-TextValue $ws.Range("D39") "0.02528"
Set-TextValue $ws.Range("E39") "1.50%"

Set-TextValue $ws.Range("D40") "0.05263"
Set-TextValue $ws.Range("E40") "1.22%"

Set-TextValue $ws.Range("D41") "0.007381"
Set-TextValue $ws.Range("E41") "-4.62%"

Set-TextValue $ws.Range("D42") "0.009616"
Set-TextValue $ws.Range("E42") "6.11%"

Set-TextValue $ws.Range("D43") "0.1404"
Set-TextValue $ws.Range("E43") "-1.53%"

Set-TextValue $ws.Range("D44") "0.002119"
Set-TextValue $ws.Range("E44") "-2.07%"

Set-TextValue $ws.Range("D45") "0.009998"
Set-TextValue $ws.Range("E45") "11.00%"

Set-TextValue $ws.Range("D46") "0.00006745"
Set-TextValue $ws.Range("E46") "1.90%"

Set-TextValue $ws.Range("D47") "0.00000000751"

Set-TextValue $ws.Range("E48") "-10.24%"

Set-TextValue $ws.Range("E49") "-7.72%"

Set-TextValue $ws.Range("D50") "0.00002102"

Set-TextValue $ws.Range("D51") "0.0002002"

